$d = $word.ActiveDocument

$replacements = @(
    @("401÷6=66, 5", "857÷3=285, 2"),
    @("883÷5=176, 3", "373÷2=186, 1"),
    @("509÷2=254, 1", "955÷3=318, 1"),
    @("176÷6=29, 2", "417÷9=46, 3"),
    @("848÷2=424, 0", "213÷9=23, 6"),
    @("764÷9=84, 8", "805÷8=100, 5"),
    @("900÷6=150, 0", "846÷7=120, 6"),
    @("849÷7=121, 2", "952÷3=317, 1"),
    @("674÷9=74, 8", "921÷4=230, 1"),
    @("290÷5=58, 0", "501÷8=62, 5"),
    @("477÷8=59, 5", "894÷3=298, 0"),
    @("377÷2=188, 1", "567÷3=189, 0"),
    @("345÷4=86, 1", "458÷8=57, 2"),
    @("528÷2=264, 0", "108÷8=13, 4"),
    @("701÷8=87, 5", "117÷2=58, 1"),
    @("449÷7=64, 1", "582÷9=64, 6"),
    @("435÷9=48, 3", "540÷4=135, 0"),
    @("497÷7=71, 0", "422÷9=46, 8"),
    @("417÷6=69, 3", "239÷3=79, 2"),
    @("161÷5=32, 1", "652÷7=93, 1"),
    @("430÷9=47, 7", "344÷2=172, 0"),
    @("782÷6=130, 2", "420÷7=60, 0"),
    @("231÷2=115, 1", "600÷4=150, 0"),
    @("303÷5=60, 3", "702÷8=87, 6"),
    @("690÷5=138, 0", "877÷4=219, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
